$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with the new block-order labels
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "living_rooms_1"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "living_rooms_2"

# Update the one-hot data rows (rows 2-7) to match the new order
$data = @(
    @(0,0,1,0,0,0),
    @(1,0,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,0,0,1,0),
    @(0,1,0,0,0,0),
    @(0,0,0,1,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
